# Hindalco price sheet update (2025-08-20)
# Inserts a new top data row (row 2) with the latest price entry, which
# shifts all existing data rows down by one. Hyperlinks that lived in
# column F are re-pointed to follow their rows, and a new hyperlink is
# added for the newly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 2, pushing everything else down ---
$ws.Rows.Item(2).Insert()

# Copy formatting from the row immediately below (the row that used to be
# row 2, now row 3) so the new row matches the existing table styling.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Populate the new row with the latest circular data ---
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 262.5
$ws.Range("E2").Value = "20.08.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf"

# --- Rebuild the hyperlinks so they track their (now shifted) rows and
#     the newly added row 2 gets its own link ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")

Write-Host "Applied Hindalco update; UsedRange:" $ws.UsedRange.Address()
